$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.404847145080566
$ws.Range("B1").Value = 3.317728281021118
$ws.Range("C1").Value = 2.498995780944824
$ws.Range("D1").Value = 2.375333786010742
$ws.Range("E1").Value = 2.444192886352539
